$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("B8").Value = "R006"
$ws.Range("B9").Value = "R007"
$ws.Range("B10").Value = "R008"

$ws.Range("C8").Value = "asdf"
$ws.Range("C9").Value = "hghhhh"
$ws.Range("C10").Value = "lllllll"

$ws.Range("N10").Value = "wayyy out here"

$ws.Range("N16").Select()
